$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.576.67"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "2.302.76"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Formula = "'537.61"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D6").Formula = "'132.12"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").Value = "2.300.63"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").Formula = "'0.101"
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D14").Formula = "'23.73"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").Value = "2.712.80"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "58.502.40"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").Value = "2.293.63"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").Formula = "'10.57"
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("E20").Value = "  -2.96%  "
$ws.Range("D21").Formula = "'315.84"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Formula = "'6.59"
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Formula = "'63.04"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Formula = "'0.170"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -2.14%  "
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").Formula = "'171.41"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("E30").Value = "  -2.93%  "
$ws.Range("D31").Value = "0.0₃0727"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("D32").Formula = "'1.11"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").Formula = "'5.84"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").Formula = "'0.385"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").Formula = "'141.08"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").Formula = "'289.55"
$ws.Range("E42").Value = "  -4.75%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").Formula = "'18.29"
$ws.Range("E47").Value = "  -2.48%  "
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("D49").Formula = "'10.96"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("E51").Value = "  -0.46%  "
